$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '67.578.46'
Set-TextValue 'E2' '  -0.13%  '
Set-TextValue 'D3' '3.781.55'
Set-TextValue 'E3' '  -0.12%  '
Set-TextValue 'E4' '  +0.11%  '
Set-TextValue 'D5' '598.15'
Set-TextValue 'E5' '  +0.41%  '
Set-TextValue 'D6' '164.49'
Set-TextValue 'E6' '  -1.33%  '
Set-TextValue 'E7' '  -0.06%  '
Set-TextValue 'E8' '  -1.07%  '
Set-TextValue 'E9' '  -0.81%  '
Set-TextValue 'E10' '  +0.28%  '
Set-TextValue 'E11' '  +0.80%  '
Set-TextValue 'D12' '0.0000248'
Set-TextValue 'E12' '  -1.71%  '
Set-TextValue 'D13' '35.48'
Set-TextValue 'E13' '  -1.43%  '
Set-TextValue 'D14' '4.414.11'
Set-TextValue 'D15' '3.778.70'
Set-TextValue 'E15' '  -0.49%  '
Set-TextValue 'D16' '67.588.25'
Set-TextValue 'E16' '  -0.08%  '
Set-TextValue 'D17' '18.31'
Set-TextValue 'E17' '  -1.27%  '
Set-TextValue 'E18' '  +1.67%  '
Set-TextValue 'D19' '7.03'
Set-TextValue 'E19' '  -0.37%  '
Set-TextValue 'D20' '460.23'
Set-TextValue 'E20' '  +0.32%  '
Set-TextValue 'D21' '9.70'
Set-TextValue 'E21' '  -2.80%  '
Set-TextValue 'D22' '0.693'
Set-TextValue 'E22' '  -0.63%  '
Set-TextValue 'D23' '0.0000146'
Set-TextValue 'E23' '  -3.87%  '
Set-TextValue 'D24' '82.50'
Set-TextValue 'E24' '  -1.00%  '
Set-TextValue 'D25' '11.97'
Set-TextValue 'E25' '  -0.57%  '
Set-TextValue 'D26' '2.09'
Set-TextValue 'E26' '  -0.67%  '
Set-TextValue 'E27' '  -0.10%  '
Set-TextValue 'D28' '9.92'
Set-TextValue 'E28' '  -0.79%  '
Set-TextValue 'D29' '3.927.89'
Set-TextValue 'E29' '  -0.18%  '
Set-TextValue 'D30' '7.42'
Set-TextValue 'E30' '  +2.96%  '
Set-TextValue 'E31' '  -6.12%  '
Set-TextValue 'E32' '  -3.27%  '
Set-TextValue 'D33' '29.04'
Set-TextValue 'E33' '  -1.72%  '
Set-TextValue 'D34' '1.00'
Set-TextValue 'E34' '  +0.30%  '
Set-TextValue 'D35' '8.93'
Set-TextValue 'E35' '  -1.29%  '
Set-TextValue 'D36' '0.0988'
Set-TextValue 'E36' '  -1.09%  '
Set-TextValue 'E37' '  +0.21%  '
Set-TextValue 'D38' '3.25'
Set-TextValue 'E38' '  -2.55%  '
Set-TextValue 'D39' '0.984'
Set-TextValue 'E39' '  -0.80%  '
Set-TextValue 'D40' '5.75'
Set-TextValue 'E40' '  -0.27%  '
Set-TextValue 'D41' '0.999'
Set-TextValue 'E41' '  -0.04%  '
Set-TextValue 'E42' '  +0.02%  '
Set-TextValue 'D43' '47.45'
Set-TextValue 'E43' '  -1.19%  '
Set-TextValue 'D44' '43.37'
Set-TextValue 'E44' '  -1.07%  '
Set-TextValue 'D45' '0.296'
Set-TextValue 'E45' '  -0.08%  '
Set-TextValue 'D46' '151.83'
Set-TextValue 'E46' '  +0.74%  '
Set-TextValue 'D47' '8.32'
Set-TextValue 'E47' '  +0.54%  '
Set-TextValue 'B48' 'ONDO'
Set-TextValue 'C48' 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 'D48' '1.35'
Set-TextValue 'E48' '  +6.68%  '
Set-TextValue 'B49' 'EnergySwap'
Set-TextValue 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '26.87'
Set-TextValue 'E49' '  +0.56%  '
Set-TextValue 'B50' 'Stacks'
Set-TextValue 'C50' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D50' '1.84'
Set-TextValue 'E50' '  +1.29%  '
Set-TextValue 'D51' '390.43'
Set-TextValue 'E51' '  +0.50%  '
